$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Step A: simple text replacements (rows 1-5) ---
$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"
$t.Cell(4,1).Range.Text = "1480"
$t.Cell(5,1).Range.Text = "0.00001"
# row 6 (0.00249) is unchanged

# --- Step B: delete old rows 7, 8, 9 (texts 0.00009 / 0.00007 / 0.00009) ---
$t.Rows.Item(7).Delete()
$t.Rows.Item(7).Delete()
$t.Rows.Item(7).Delete()

# After the deletes: old row10 (0.00010) is now row 7 (text unchanged),
# old row11 (0.00012) is now row 8, old row12 (0.06677) is now row 9.

# --- Step C: edit the shifted rows 8 and 9 ---
$t.Cell(8,1).Range.Text = "0.00004"
$t.Cell(9,1).Range.Text = "0.00017"

# --- Step D: insert 3 fresh rows after row 9 (before row 10) ---
$t.Rows.Add($t.Rows.Item(10)) | Out-Null
$t.Cell(10,1).Range.Text = "0.00019"

$t.Rows.Add($t.Rows.Item(11)) | Out-Null
$t.Cell(11,1).Range.Text = "0.00020"

$t.Rows.Add($t.Rows.Item(12)) | Out-Null
$t.Cell(12,1).Range.Text = "0.17709"

# Table is back to 46 rows, so rows 13-43 keep their original indices/content.

# --- Step E: collapse the big tab-separated rows (now still at rows 44-46) ---
$t.Cell(44,1).Range.Text = "99.88"
$t.Cell(45,1).Range.Text = "0.18"
$t.Cell(46,1).Range.Text = "146"
